# Auto-generated Excel COM-interop script to apply the cryptos.xlsx price/volume update
# (commit: "Updated cryptos list on Mon Dec  4 23:38:21 UTC 2023 with GitHub Actions")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing it to stay a TEXT cell, even
# when the text looks like a plain number (e.g. "59.00" or "0.998"). Excel
# normally auto-converts such strings to numbers on a bare .Value assignment,
# which would silently drop significant trailing zeros / dot-grouping. Briefly
# forcing a Text number format makes the assignment keep the literal text, and
# ClearFormats() immediately after removes that temporary formatting again so
# the cell style is left exactly as it was (default/general).
function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") "41.905.16"
$ws.Range("E2").Value = "  +4.81%  "

# Row 3
Set-TextValue $ws.Range("D3") "2.237.36"
$ws.Range("E3").Value = "  +2.19%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
Set-TextValue $ws.Range("D5") "232.89"
$ws.Range("E5").Value = "  +2.35%  "

# Row 6
$ws.Range("E6").Value = "  -0.89%  "

# Row 7
Set-TextValue $ws.Range("D7") "61.58"
$ws.Range("E7").Value = "  -2.62%  "

# Row 8
$ws.Range("E8").Value = "  +0.06%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.408"
$ws.Range("E9").Value = "  +3.27%  "

# Row 10
Set-TextValue $ws.Range("D10") "59.00"
$ws.Range("E10").Value = "  +1.30%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.0920"
$ws.Range("E11").Value = "  +7.00%  "

# Row 12
$ws.Range("E12").Value = "  +0.66%  "

# Row 13
Set-TextValue $ws.Range("D13") "2.565.00"
$ws.Range("E13").Value = "  +1.87%  "

# Row 14
Set-TextValue $ws.Range("D14") "15.72"
$ws.Range("E14").Value = "  -0.19%  "

# Row 15
Set-TextValue $ws.Range("D15") "22.40"
$ws.Range("E15").Value = "  +1.94%  "

# Row 16
Set-TextValue $ws.Range("D16") "0.807"
$ws.Range("E16").Value = "  -0.59%  "

# Row 17
Set-TextValue $ws.Range("D17") "5.63"
$ws.Range("E17").Value = "  +1.60%  "

# Row 18
Set-TextValue $ws.Range("D18") "2.240.20"
$ws.Range("E18").Value = "  +1.78%  "

# Row 19
Set-TextValue $ws.Range("D19") "41.791.11"
$ws.Range("E19").Value = "  +4.83%  "

# Row 20
Set-TextValue $ws.Range("D20") "0.0₃0915"
$ws.Range("E20").Value = "  +0.99%  "

# Row 21
Set-TextValue $ws.Range("D21") "72.67"
$ws.Range("E21").Value = "  +0.81%  "

# Row 22
Set-TextValue $ws.Range("D22") "6.06"
$ws.Range("E22").Value = "  +0.28%  "

# Row 23
Set-TextValue $ws.Range("D23") "253.91"
$ws.Range("E23").Value = "  +9.34%  "

# Row 24
$ws.Range("E24").Value = "  +0.04%  "

# Row 25
Set-TextValue $ws.Range("D25") "2.40"
$ws.Range("E25").Value = "  +2.48%  "

# Row 26
Set-TextValue $ws.Range("D26") "2.37"
$ws.Range("E26").Value = "  -0.71%  "

# Row 27
Set-TextValue $ws.Range("D27") "9.75"
$ws.Range("E27").Value = "  +1.17%  "

# Row 28
Set-TextValue $ws.Range("D28") "0.145"
$ws.Range("E28").Value = "  +3.94%  "

# Row 29
Set-TextValue $ws.Range("D29") "169.51"
$ws.Range("E29").Value = "  -1.06%  "

# Row 30
Set-TextValue $ws.Range("D30") "20.11"
$ws.Range("E30").Value = "  +0.32%  "

# Row 31
Set-TextValue $ws.Range("D31") "1.43"
$ws.Range("E31").Value = "  -1.33%  "

# Row 32
Set-TextValue $ws.Range("D32") "2.73"
$ws.Range("E32").Value = "  -0.43%  "

# Row 33
$ws.Range("E33").Value = "  -0.15%  "

# Row 34
Set-TextValue $ws.Range("D34") "5.08"
$ws.Range("E34").Value = "  +8.15%  "

# Row 35
Set-TextValue $ws.Range("D35") "4.70"
$ws.Range("E35").Value = "  +3.27%  "

# Row 36
Set-TextValue $ws.Range("D36") "0.0640"
$ws.Range("E36").Value = "  +2.75%  "

# Row 37
Set-TextValue $ws.Range("D37") "6.67"
$ws.Range("E37").Value = "  -4.70%  "

# Row 38
Set-TextValue $ws.Range("D38") "3.75"
$ws.Range("E38").Value = "  -3.00%  "

# Row 39
Set-TextValue $ws.Range("D39") "2.37"
$ws.Range("E39").Value = "  -2.40%  "

# Row 40
Set-TextValue $ws.Range("D40") "0.000261"
$ws.Range("E40").Value = "  +34.21%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.998"
$ws.Range("E41").Value = "  -0.27%  "

# Row 42
Set-TextValue $ws.Range("D42") "0.0242"
$ws.Range("E42").Value = "  +5.83%  "

# Row 43
$ws.Range("B43").Value = "FTXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws.Range("D43") "4.71"
$ws.Range("E43").Value = "  -6.39%  "

# Row 44
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D44") "8.62"
$ws.Range("E44").Value = "  +4.44%  "

# Row 45
$ws.Range("E45").Value = "  +1.13%  "

# Row 46
Set-TextValue $ws.Range("D46") "99.79"
$ws.Range("E46").Value = "  -3.16%  "

# Row 47
Set-TextValue $ws.Range("D47") "0.0960"
$ws.Range("E47").Value = "  +3.27%  "

# Row 48
Set-TextValue $ws.Range("D48") "1.485.87"
$ws.Range("E48").Value = "  -1.92%  "

# Row 49
Set-TextValue $ws.Range("D49") "16.61"
$ws.Range("E49").Value = "  -4.44%  "

# Row 50
$ws.Range("E50").Value = "  +0.15%  "

# Row 51
Set-TextValue $ws.Range("D51") "52.83"
$ws.Range("E51").Value = "  +6.31%  "
